# Add a new "waiter" row (row 2) to the "User Data" sheet:
#   A=User ID  B=Email  C=First Name  D=Last Name  E=Staff ID
#   F=User Type  G=Hours Worked  H=Total Hours  I=Active Status
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User Data")

$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "ww"
$ws.Range("C2").Value = "ww"
$ws.Range("D2").Value = "ww"
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = "Waiter"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 20
$ws.Range("I2").Value = $true

# Resize columns to fit the new, narrower content (mirrors the author
# double-click-to-autofit on columns B, E-H after entering the new row).
$ws.Range("B:B").ColumnWidth = 5.0
$ws.Range("E:E").ColumnWidth = 6.666666666666667
$ws.Range("F:F").ColumnWidth = 9.166666666666666
$ws.Range("G:G").ColumnWidth = 13.166666666666666
$ws.Range("H:H").ColumnWidth = 10.5

# Leave the selection on B2, matching the saved workbook's cursor position.
$null = $ws.Range("B2").Select()
